$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column AI (base_damage_mod_bonus) values to 0 for rows 3 through 60
for ($r = 3; $r -le 60; $r++) {
    $ws.Range("AI$r").Value = 0
}

# Update the view: scroll position and selection
$ws.Application.ActiveWindow.ScrollRow = 39
$ws.Application.ActiveWindow.ScrollColumn = 16
$ws.Range("AI2:AI60").Select()
